$d = $word.ActiveDocument

# The paragraph currently holds the id tag split across three runs:
#   <id>  |  p001v_1  |  </id>
# with the surrounding "<id>" / "</id>" runs in Courier New / color 7f6000
# and the "p001v_1" run in the default body formatting. The edit merges
# them into a single run (keeping the first run's Courier-New/7f6000
# formatting) containing the full "<id>p001v_1</id>" text.

$rng = $d.Content
$found = $rng.Find.Execute("<id>p001v_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $finalText = $rng.Text

    # Re-assigning the exact same text is a no-op, so nudge the range
    # through a distinct value first - this collapses the three runs into
    # one (inheriting the leading run's character formatting), then we set
    # the final text on that single, now-unified run.
    $rng.Text = $finalText + "#"
    $rng.Text = $finalText
}
